$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Unidades Pedido" (L) and "Diferencia Stock" (M) values for the
# affected article rows.
$ws.Range("L2").Value = 2
$ws.Range("M2").Value = 1

$ws.Range("L3").Value = 7
$ws.Range("M3").Value = 1

$ws.Range("L4").Value = 6
$ws.Range("M4").Value = 1

$ws.Range("L6").Value = 2
$ws.Range("M6").Value = 1

$ws.Range("L7").Value = 5
$ws.Range("M7").Value = 1

$ws.Range("L16").Value = 2
$ws.Range("M16").Value = 1

$ws.Range("L17").Value = 6
$ws.Range("M17").Value = 1

$ws.Range("L18").Value = 49
$ws.Range("M18").Value = 6

$ws.Range("L19").Value = 4
$ws.Range("M19").Value = 1

$ws.Range("L21").Value = 5
$ws.Range("M21").Value = -3

$ws.Range("L22").Value = 5
$ws.Range("M22").Value = 1

$ws.Range("L24").Value = 5
$ws.Range("M24").Value = 1

$ws.Range("L27").Value = 11
$ws.Range("M27").Value = 1

$ws.Range("L34").Value = 6
$ws.Range("M34").Value = 1

$ws.Range("L35").Value = 6
$ws.Range("M35").Value = 1

# Update the summary totals.
$ws.Range("C40").Value = 174
$ws.Range("C51").Value = 16
